# "Update countries & provincias Spain"
# Refresh the COVID-19 stats table on sheet "Pais": update the timestamp
# banner and the per-country Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes figures (columns
# B:H). Because the sheet is kept sorted by "Casos totales" (column B)
# descending, a handful of neighbouring countries swap row positions once
# their totals are refreshed (Costa Rica/Camerun/Venezuela around row 69-71,
# Maldivas/Nicaragua around row 107-108, Angola/Burkina Faso/Republica de
# Chipre around row 146-148) - those rows get a new country name (column A)
# as well as new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Banner timestamp (row 1)
$ws.Range("A1").Value = 'Datos actualizados a 30 de Julio de 2020 a las 22:28'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4613251
$ws.Range("C4").Value = 45214
$ws.Range("D4").Value = 2257949
$ws.Range("E4").Value = 2200521
$ws.Range("G4").Value = 941
$ws.Range("H4").Value = 154781

# Israel (row 36)
$ws.Range("B36").Value = 70036
$ws.Range("C36").Value = 1737
$ws.Range("E36").Value = 34020
$ws.Range("G36").Value = 9
$ws.Range("H36").Value = 500

# Moldavia (row 63)
$ws.Range("E63").Value = 6136
$ws.Range("G63").Value = 12
$ws.Range("H63").Value = 771

# Costa Rica moves up to row 69
$ws.Range("A69").Value = 'Costa Rica'
$ws.Range("B69").Value = 17290
$ws.Range("C69").Value = 490
$ws.Range("D69").Value = 4280
$ws.Range("E69").Value = 12870
$ws.Range("G69").Value = 7
$ws.Range("H69").Value = 140

# Camerun shifts down to row 70
$ws.Range("A70").Value = 'Camerun'
$ws.Range("B70").Value = 17255
$ws.Range("D70").Value = 15320
$ws.Range("E70").Value = 1544
$ws.Range("H70").Value = 391

# Venezuela shifts down to row 71
$ws.Range("A71").Value = 'Venezuela'
$ws.Range("B71").Value = 17158
$ws.Range("D71").Value = 10421
$ws.Range("E71").Value = 6583
$ws.Range("H71").Value = 154

# Luxemburgo (row 95)
$ws.Range("B95").Value = 6616
$ws.Range("C95").Value = 83
$ws.Range("D95").Value = 5027
$ws.Range("E95").Value = 1475

# Mauritania (row 96)
$ws.Range("B96").Value = 6295
$ws.Range("C96").Value = 22
$ws.Range("D96").Value = 4889
$ws.Range("E96").Value = 1249
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 157

# Maldivas moves up to row 107
$ws.Range("A107").Value = 'Maldivas'
$ws.Range("B107").Value = 3719
$ws.Range("C107").Value = 152
$ws.Range("D107").Value = 2568
$ws.Range("E107").Value = 1135
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 16

# Nicaragua shifts down to row 108
$ws.Range("A108").Value = 'Nicaragua'
$ws.Range("B108").Value = 3672
$ws.Range("D108").Value = 2492
$ws.Range("E108").Value = 1064
$ws.Range("H108").Value = 116

# Yemen (row 135)
$ws.Range("B135").Value = 1726
$ws.Range("C135").Value = 15
$ws.Range("D135").Value = 856
$ws.Range("E135").Value = 383
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = 487

# Angola moves up to row 146
$ws.Range("A146").Value = 'Angola'
$ws.Range("B146").Value = 1109
$ws.Range("C146").Value = 31
$ws.Range("D146").Value = 395
$ws.Range("E146").Value = 663
$ws.Range("G146").Value = 3
$ws.Range("H146").Value = 51

# Burkina Faso shifts down to row 147
$ws.Range("A147").Value = 'Burkina Faso'
$ws.Range("B147").Value = 1106
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 935
$ws.Range("E147").Value = 118
$ws.Range("H147").Value = 53

# Republica de Chipre shifts down to row 148
$ws.Range("A148").Value = 'Republica de Chipre'
$ws.Range("B148").Value = 1090
$ws.Range("C148").Value = 10
$ws.Range("D148").Value = 852
$ws.Range("E148").Value = 219
$ws.Range("H148").Value = 19

# Republica del Chad (row 149)
$ws.Range("B149").Value = 935
$ws.Range("C149").Value = 9
$ws.Range("E149").Value = 47

# Togo (row 151)
$ws.Range("B151").Value = 908
$ws.Range("C151").Value = 12
$ws.Range("D151").Value = 626
$ws.Range("E151").Value = 264

# Santo Tome y Principe (row 152)
$ws.Range("B152").Value = 870
$ws.Range("C152").Value = 2
$ws.Range("D152").Value = 771
$ws.Range("E152").Value = 84

# Trinidad yTobago (row 180)
$ws.Range("B180").Value = 157
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 130
$ws.Range("E180").Value = 19
